$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.88
$ws.Range("A14").Value = -21.695
$ws.Range("B15").Value = 5.792
$ws.Range("A16").Value = -22.181
$ws.Range("A21").Value = -21.204
$ws.Range("B21").Value = 7.292
$ws.Range("B22").Value = 7.290999999999999
$ws.Range("A23").Value = -20.812
$ws.Range("B24").Value = 5.513
$ws.Range("A25").Value = -21.589
$ws.Range("A26").Value = -21.083
$ws.Range("B27").Value = 6.44
$ws.Range("B28").Value = 5.583000000000001
$ws.Range("A29").Value = -21.52
$ws.Range("B36").Value = 6.356
$ws.Range("B39").Value = 6.609
$ws.Range("A40").Value = -20.706
$ws.Range("B45").Value = 5.832
$ws.Range("B48").Value = 5.544
$ws.Range("B49").Value = 6.470000000000001
$ws.Range("B52").Value = 5.968000000000001
$ws.Range("A53").Value = -20.623
$ws.Range("B53").Value = 7.706999999999999
$ws.Range("B54").Value = 5.077000000000001
$ws.Range("A57").Value = -22.053
$ws.Range("B57").Value = 5.720000000000001
$ws.Range("A59").Value = -22.437
$ws.Range("A65").Value = -21.529
$ws.Range("A69").Value = -21.476
$ws.Range("B70").Value = 5.029999999999999
$ws.Range("B71").Value = 5.457000000000001
$ws.Range("A79").Value = -21.285
$ws.Range("A83").Value = -21.979
$ws.Range("B86").Value = 4.937
$ws.Range("B87").Value = 4.901999999999999
$ws.Range("B89").Value = 4.759
$ws.Range("A91").Value = -21.033
$ws.Range("A93").Value = -21.533
$ws.Range("A100").Value = -22.368
$ws.Range("B101").Value = 5.279000000000001
$ws.Range("A103").Value = -22.055
